$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (shifts LOC_2019..LOC_2023 down by one row)
$ws.Rows.Item(6).Insert()

# Every value on this sheet (including dates and numbers) is stored as
# plain text, so force Text format before assigning the new row's values
# to stop Excel from auto-converting them to real dates/numbers.
$ws.Range("A6:J6").NumberFormat = "@"

# Fill in the new LOC_2018 row
$ws.Range("A6").Value = "LOC_2018"
$ws.Range("B6").Value = "2018-05-01"
$ws.Range("C6").Value = "2018-05-21"
$ws.Range("D6").Value = "249.24"
$ws.Range("E6").Value = "449.07"
$ws.Range("F6").Value = "0.966729365689967"
$ws.Range("G6").Value = "1.67136058738539e-15"
$ws.Range("H6").Value = "0.000124804743867283"
$ws.Range("I6").Value = "-190122.828179172"
$ws.Range("J6").Value = "full_ice_to_no_ice"

# Restore the default (unstyled) look for the new row to match the rest
# of the data rows, which carry no explicit cell style.
$ws.Range("A6:J6").Style = "Normal"
